$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 579; this shifts the existing rows 579-630
# down to 580-631 (preserving all their values/formatting) and keeps the
# sheet dimension in sync (A1:R630 -> A1:R631).
$ws.Rows(579).Insert()

# Populate the newly inserted row 579 with the new data point.
$ws.Cells.Item(579, 1).Value  = 3
$ws.Cells.Item(579, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(579, 3).Value  = "Coquimbo"
$ws.Cells.Item(579, 4).Value  = 45166
$ws.Cells.Item(579, 5).Value  = 5
$ws.Cells.Item(579, 6).Value  = 100112040
$ws.Cells.Item(579, 7).Value  = "Cilantro"
$ws.Cells.Item(579, 8).Value  = "Sin especificar"
$ws.Cells.Item(579, 9).Value  = "Primera"
$ws.Cells.Item(579, 10).Value = 120
$ws.Cells.Item(579, 11).Value = 4000
$ws.Cells.Item(579, 12).Value = 4000
$ws.Cells.Item(579, 13).Value = 4000
$ws.Cells.Item(579, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(579, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(579, 16).Value = 1333
$ws.Cells.Item(579, 17).Value = 3
$ws.Cells.Item(579, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date number-format as the rest
# of column D (style index 2 in the original workbook).
$ws.Cells.Item(579, 4).NumberFormat = $ws.Cells.Item(580, 4).NumberFormat
